# Auto-generated edit script: updates currentAveragePrice* / Leve*Profit* columns
# per the scheduled-runner data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 324.25
$ws.Range("I4").Value = 209.6
$ws.Range("K4").Value = 209.6
$ws.Range("M4").Value = -95.59999999999999
$ws.Range("H5").Value = 132.33333
$ws.Range("I5").Value = 132.33333
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 132.33333
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = $null
$ws.Range("N5").Value = -17.33332999999999
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = $null
$ws.Range("H17").Value = 1778.3889
$ws.Range("J17").Value = 2143.9048
$ws.Range("L17").Value = 6431.714399999999
$ws.Range("N17").Value = -6767.714399999999
$ws.Range("H18").Value = 843
$ws.Range("I18").Value = 873.2
$ws.Range("J18").Value = 692
$ws.Range("K18").Value = 873.2
$ws.Range("L18").Value = 692
$ws.Range("M18").Value = -589.2
$ws.Range("N18").Value = -1260
$ws.Range("H38").Value = 344.35715
$ws.Range("I38").Value = 178.53847
$ws.Range("J38").Value = 2500
$ws.Range("K38").Value = 535.61541
$ws.Range("L38").Value = 7500
$ws.Range("M38").Value = -163.61541
$ws.Range("N38").Value = -8244
$ws.Range("H39").Value = 139.58824
$ws.Range("I39").Value = 33.785713
$ws.Range("K39").Value = 101.357139
$ws.Range("M39").Value = 194.642861
$ws.Range("H43").Value = 2000
$ws.Range("I43").Value = 2000
$ws.Range("K43").Value = 2000
$ws.Range("M43").Value = -1931
$ws.Range("H55").Value = 100
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 100
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = $null
$ws.Range("N55").Value = 114
$ws.Range("H76").Value = 2527.5715
$ws.Range("I76").Value = 2782.1667
$ws.Range("J76").Value = 1000
$ws.Range("K76").Value = 2782.1667
$ws.Range("L76").Value = 1000
$ws.Range("M76").Value = -2467.1667
$ws.Range("N76").Value = -1630
$ws.Range("H79").Value = 2527.5715
$ws.Range("I79").Value = 2782.1667
$ws.Range("J79").Value = 1000
$ws.Range("K79").Value = 2782.1667
$ws.Range("L79").Value = 1000
$ws.Range("M79").Value = -1690.1667
$ws.Range("N79").Value = -3184
$ws.Range("H125").Value = 4333.1665
$ws.Range("I125").Value = 3999.75
$ws.Range("K125").Value = 35997.75
$ws.Range("M125").Value = -33537.75
$ws.Range("H135").Value = 2398
$ws.Range("I135").Value = 2462.25
$ws.Range("J135").Value = 2333.75
$ws.Range("K135").Value = 22160.25
$ws.Range("L135").Value = 21003.75
$ws.Range("M135").Value = -19625.25
$ws.Range("N135").Value = -26073.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7219
$ws.Range("I61").Value = 7315.6665
$ws.Range("K61").Value = 7315.6665
$ws.Range("M61").Value = -7103.6665
$ws.Range("H74").Value = 4064.4
$ws.Range("I74").Value = 2733.111
$ws.Range("K74").Value = 2733.111
$ws.Range("M74").Value = -1859.111
$ws.Range("H77").Value = 4064.4
$ws.Range("I77").Value = 2733.111
$ws.Range("K77").Value = 13665.555
$ws.Range("M77").Value = -9297.555
$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51802
$ws.Range("H102").Value = 6748.5835
$ws.Range("I102").Value = 5247.875
$ws.Range("K102").Value = 5247.875
$ws.Range("M102").Value = -3625.875
$ws.Range("H110").Value = 4195.3335
$ws.Range("I110").Value = 3960
$ws.Range("J110").Value = 4666
$ws.Range("K110").Value = 3960
$ws.Range("L110").Value = 4666
$ws.Range("M110").Value = -1915
$ws.Range("N110").Value = -8756
$ws.Range("H132").Value = 2784.1538
$ws.Range("I132").Value = 1719.5
$ws.Range("K132").Value = 5158.5
$ws.Range("M132").Value = -2628.5
$ws.Range("H136").Value = 7219
$ws.Range("I136").Value = 7315.6665
$ws.Range("K136").Value = 21946.9995
$ws.Range("M136").Value = -19396.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 14340.75
$ws.Range("J97").Value = 40000
$ws.Range("L97").Value = 40000
$ws.Range("N97").Value = -41982

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4173.35
$ws.Range("I132").Value = 3558.6
$ws.Range("J132").Value = 6017.6
$ws.Range("K132").Value = 10675.8
$ws.Range("L132").Value = 18052.8
$ws.Range("M132").Value = -8145.799999999999
$ws.Range("N132").Value = -23112.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 3533.6155
$ws.Range("I92").Value = 3897.6667
$ws.Range("J92").Value = 3221.5715
$ws.Range("K92").Value = 11693.0001
$ws.Range("L92").Value = 9664.7145
$ws.Range("M92").Value = -10445.0001
$ws.Range("N92").Value = -12160.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 10250
$ws.Range("I21").Value = 3500
$ws.Range("K21").Value = 3500
$ws.Range("M21").Value = -3327
$ws.Range("H30").Value = 10250
$ws.Range("I30").Value = 3500
$ws.Range("K30").Value = 3500
$ws.Range("M30").Value = -3395
$ws.Range("H113").Value = 7407.25
$ws.Range("I113").Value = 4832.6665
$ws.Range("J113").Value = 8265.444
$ws.Range("K113").Value = 4832.6665
$ws.Range("L113").Value = 8265.444
$ws.Range("M113").Value = -2662.6665
$ws.Range("N113").Value = -12605.444

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 7007.4443
$ws.Range("J25").Value = 9508
$ws.Range("L25").Value = 9508
$ws.Range("N25").Value = -9968
$ws.Range("H40").Value = 6777.2104
$ws.Range("J40").Value = 8956
$ws.Range("L40").Value = 8956
$ws.Range("N40").Value = -9228
$ws.Range("H82").Value = 3143.3125
$ws.Range("J82").Value = 3649.7
$ws.Range("L82").Value = 3649.7
$ws.Range("N82").Value = -4371.7
$ws.Range("H85").Value = 3143.3125
$ws.Range("J85").Value = 3649.7
$ws.Range("L85").Value = 3649.7
$ws.Range("N85").Value = -6145.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2214.1428
$ws.Range("I2").Value = 2214.1428
$ws.Range("K2").Value = 2214.1428
$ws.Range("M2").Value = -2102.1428
$ws.Range("H4").Value = 10509.5
$ws.Range("I4").Value = 14899.286
$ws.Range("J4").Value = 266.66666
$ws.Range("K4").Value = 14899.286
$ws.Range("L4").Value = 266.66666
$ws.Range("M4").Value = -14786.286
$ws.Range("N4").Value = -492.66666
$ws.Range("H129").Value = 550000
$ws.Range("J129").Value = 550000
$ws.Range("L129").Value = 550000
$ws.Range("N129").Value = -560000

Write-Host "Applied 172 cell updates across 8 sheets"
